$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grievanceDetails")
$ws.Activate()

# Replace the stale ward/location value with the new grievance location.
$ws.Range("E2").Value = "Aavanthi Nagar"

# Author also left the active selection on F6 when they saved.
$ws.Range("F6").Select()
